# Atualizando tabela e adicionando status do time no campeonato
# Preenche as colunas D (Score_m) e F (Score_v) para as linhas 182-201

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scores = @{
    182 = @(0, 1)
    183 = @(1, 1)
    184 = @(3, 1)
    185 = @(1, 0)
    186 = @(3, 1)
    187 = @(1, 1)
    188 = @(1, 0)
    189 = @(2, 1)
    190 = @(2, 0)
    191 = @(4, 0)
    192 = @(0, 2)
    193 = @(4, 0)
    194 = @(2, 3)
    195 = @(0, 1)
    196 = @(1, 1)
    197 = @(1, 1)
    198 = @(2, 1)
    199 = @(1, 0)
    200 = @(0, 0)
    201 = @(3, 1)
}

foreach ($row in $scores.Keys) {
    $pair = $scores[$row]

    # Use the formatting already present in column A (style index 2) so the
    # newly created cells match the rest of the table instead of falling
    # back to Excel's default style.
    $srcCell = $ws.Cells.Item($row, 1)

    $dCell = $ws.Cells.Item($row, 4)
    $srcCell.Copy($dCell)
    $dCell.Value = $pair[0]

    $fCell = $ws.Cells.Item($row, 6)
    $srcCell.Copy($fCell)
    $fCell.Value = $pair[1]
}
